# Add a "Hydropathy index" column (O) to the amino_acids sheet, with
# per-row Kyte-Doolittle values, a header comment describing the column,
# and widen the sheet's AutoFilter / _FilterDatabase range to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("amino_acids")

# --- New header cell (O2), formatted the same way as the other headers ---
$ws.Range("O2").Value = "Hydropathy index"
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- New data column (Kyte & Doolittle hydropathy index) ---
$values = @{
    3  = -4.5
    4  = -3.9
    5  = -3.2
    6  = -3.5
    7  = -3.5
    8  = -0.8
    9  = -0.7
    10 = -3.5
    11 = -3.5
    12 = 2.5
    14 = -0.4
    15 = -1.6
    16 = 1.8
    17 = 4.5
    18 = 3.8
    19 = 1.9
    20 = 2.8
    21 = -1.3
    22 = -0.9
    23 = 4.2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 15).Value = $values[$row]
}

# --- Explanatory comment on the new header cell ---
$commentText = "Jeromel, Luka:" + "`nMeasure of hydrophobic effect. High positive value means strong hidrophobicity."
$ws.Range("O2").AddComment($commentText)

# --- Re-point the AutoFilter (and the hidden _FilterDatabase name it
#     maintains) at the widened A2:O2 range; drop the old sort state ---
$ws.AutoFilterMode = $false
$ws.Range("A2:O2").AutoFilter()
$wb.Names.Item("amino_acids!_FilterDatabase").RefersTo = "=amino_acids!`$A`$2:`$O`$2"

# --- Restore the cursor to where the author left it ---
$ws.Range("D28").Select()
